# issue #5: stock data output to json file
# Insert a new "property_category" column into the "股票" (Stock) worksheet,
# between the "total" and "date" columns, and populate it with "stock".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Shift existing columns H:J (date, legislator_name, legislator_id) one column
# to the right (to I:K), opening up column H for the new field.
$ws.Range("H1:H2").Insert(-4161)

# New header + value for the inserted "property_category" column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
